# Generate Report for Handoff
# Refresh the "Latest Handoff Datetime" / "Latest Handoff Date" values for every
# file that is currently pending handoff ("Ready for handoff") or stuck
# ("Handback transform failed"), across the Overview sheet and each locale
# sheet, to reflect a freshly generated handoff report.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: column D = "Latest Handoff Date" ---
$overview = $wb.Worksheets.Item("Overview")
$overviewDate = "2016-22-17 14:22:46"
$overviewRows = @(7, 10, 11, 12, 13, 14, 15, 16)
foreach ($r in $overviewRows) {
    $overview.Cells.Item($r, 4).Value = $overviewDate
}

# --- zh-cn sheet: column E = "Latest Handoff Datetime" ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcnDate = "2016-03-17 14:22:42"
$localeRows = @(7, 10, 11, 12, 13, 14, 15, 16)
foreach ($r in $localeRows) {
    $zhcn.Cells.Item($r, 5).Value = $zhcnDate
}

# --- de-de sheet: column E = "Latest Handoff Datetime" ---
$dede = $wb.Worksheets.Item("de-de")
$dedeDate = "2016-03-17 14:22:46"
foreach ($r in $localeRows) {
    $dede.Cells.Item($r, 5).Value = $dedeDate
}
